$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 1063.25
$ws.Range("I29").Value = 916.6667
$ws.Range("J29").Value = 1151.2
$ws.Range("K29").Value = 2750.0001
$ws.Range("L29").Value = 3453.6
$ws.Range("M29").Value = -2469.0001
$ws.Range("N29").Value = -4015.6

# Row 51
$ws.Range("H51").Value = 5368.1
$ws.Range("I51").Value = 3565.1667
$ws.Range("K51").Value = 3565.1667
$ws.Range("M51").Value = -3081.1667

# Row 64
$ws.Range("H64").Value = 1500
$ws.Range("I64").Value = 1500
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 1500
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -1252
$ws.Range("N64").ClearContents()

# Row 67
$ws.Range("H67").Value = 1500
$ws.Range("I67").Value = 1500
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 1500
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -642
$ws.Range("N67").ClearContents()

# Row 112
$ws.Range("H112").Value = 1312.4108
$ws.Range("J112").Value = 1312.4108
$ws.Range("L112").Value = 3937.2324
$ws.Range("N112").Value = -6153.232400000001

# Row 116
$ws.Range("H116").Value = 7286.8945
$ws.Range("I116").Value = 2439.8333
$ws.Range("J116").Value = 9524
$ws.Range("K116").Value = 2439.8333
$ws.Range("L116").Value = 9524
$ws.Range("M116").Value = 1002.1667
$ws.Range("N116").Value = -16408

# Row 132
$ws.Range("H132").Value = 34833516
$ws.Range("I132").Value = 45461024
$ws.Range("J132").Value = 1432788.1
$ws.Range("K132").Value = 136383072
$ws.Range("L132").Value = 4298364.300000001
$ws.Range("M132").Value = -136380542
$ws.Range("N132").Value = -4303424.300000001

# Row 137
$ws.Range("H137").Value = 1589542
$ws.Range("I137").Value = 1765165
$ws.Range("J137").Value = 8934.333000000001
$ws.Range("K137").Value = 5295495
$ws.Range("L137").Value = 26802.999
$ws.Range("M137").Value = -5292945
$ws.Range("N137").Value = -31902.999

$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 8660957
$ws.Range("I63").Value = 17315914
$ws.Range("J63").Value = 5999.875
$ws.Range("K63").Value = 17315914
$ws.Range("L63").Value = 5999.875
$ws.Range("M63").Value = -17315228
$ws.Range("N63").Value = -7371.875

# Row 66
$ws.Range("H66").Value = 8660957
$ws.Range("I66").Value = 17315914
$ws.Range("J66").Value = 5999.875
$ws.Range("K66").Value = 86579570
$ws.Range("L66").Value = 29999.375
$ws.Range("M66").Value = -86576138
$ws.Range("N66").Value = -36863.375

# Row 97
$ws.Range("H97").Value = 1301.4546
$ws.Range("I97").Value = 1244.3334
$ws.Range("J97").Value = 1370
$ws.Range("K97").Value = 1244.3334
$ws.Range("L97").Value = 1370
$ws.Range("M97").Value = -748.3334
$ws.Range("N97").Value = -2362

# Row 122
$ws.Range("H122").Value = 3843.3958
$ws.Range("I122").Value = 3293.7026
$ws.Range("J122").Value = 5692.364
$ws.Range("K122").Value = 9881.1078
$ws.Range("L122").Value = 17077.092
$ws.Range("M122").Value = -7431.1078
$ws.Range("N122").Value = -21977.092

# Row 132
$ws.Range("H132").Value = 2641.611
$ws.Range("I132").Value = 1095.5454
$ws.Range("J132").Value = 5071.143
$ws.Range("K132").Value = 3286.6362
$ws.Range("L132").Value = 15213.429
$ws.Range("M132").Value = -756.6361999999999
$ws.Range("N132").Value = -20273.429

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 1278.3529
$ws.Range("I107").Value = 930.8570999999999
$ws.Range("K107").Value = 930.8570999999999
$ws.Range("M107").Value = 989.1429000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 763.96
$ws.Range("I22").Value = 437.4375
$ws.Range("J22").Value = 1344.4445
$ws.Range("K22").Value = 437.4375
$ws.Range("L22").Value = 1344.4445
$ws.Range("M22").Value = -87.4375
$ws.Range("N22").Value = -2044.4445

# Row 58
$ws.Range("H58").Value = 1976.0714
$ws.Range("I58").Value = 1493.7037
$ws.Range("J58").Value = 15000
$ws.Range("K58").Value = 1493.7037
$ws.Range("L58").Value = 15000
$ws.Range("M58").Value = -1290.7037
$ws.Range("N58").Value = -15406

# Row 81
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# Row 84
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# Row 132
$ws.Range("H132").Value = 5800.7144
$ws.Range("I132").Value = 2674
$ws.Range("J132").Value = 8145.75
$ws.Range("K132").Value = 8022
$ws.Range("L132").Value = 24437.25
$ws.Range("M132").Value = -5492
$ws.Range("N132").Value = -29497.25

# Row 133
$ws.Range("H133").Value = 39211.555
$ws.Range("J133").Value = 41576
$ws.Range("L133").Value = 41576
$ws.Range("N133").Value = -46636

# Row 134
$ws.Range("H134").Value = 11009.077
$ws.Range("I134").Value = 16314.857
$ws.Range("J134").Value = 4819
$ws.Range("K134").Value = 48944.571
$ws.Range("L134").Value = 14457
$ws.Range("M134").Value = -46409.571
$ws.Range("N134").Value = -19527

# Row 136
$ws.Range("H136").Value = 1976.0714
$ws.Range("I136").Value = 1493.7037
$ws.Range("J136").Value = 15000
$ws.Range("K136").Value = 4481.1111
$ws.Range("L136").Value = 45000
$ws.Range("M136").Value = -1931.1111
$ws.Range("N136").Value = -50100

$ws = $wb.Worksheets.Item("CUL")
# Row 109
$ws.Range("H109").Value = 1452.4286
$ws.Range("I109").Value = 986.4
$ws.Range("J109").Value = 2617.5
$ws.Range("K109").Value = 2959.2
$ws.Range("L109").Value = 7852.5
$ws.Range("M109").Value = -1919.2
$ws.Range("N109").Value = -9932.5

# Row 112
$ws.Range("H112").Value = 1900
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 2800
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 8400
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -10616

# Row 113
$ws.Range("H113").Value = 5952977
$ws.Range("I113").Value = 618.38464
$ws.Range("J113").Value = 15625560
$ws.Range("K113").Value = 1855.15392
$ws.Range("L113").Value = 46876680
$ws.Range("M113").Value = 314.84608
$ws.Range("N113").Value = -46881020

# Row 129
$ws.Range("H129").Value = 2128.64
$ws.Range("J129").Value = 1729.7778
$ws.Range("L129").Value = 5189.3334
$ws.Range("N129").Value = -15189.3334

# Row 132
$ws.Range("H132").Value = 2020.6757
$ws.Range("I132").Value = 865.6667
$ws.Range("J132").Value = 2808.182
$ws.Range("K132").Value = 7791.0003
$ws.Range("L132").Value = 25273.638
$ws.Range("M132").Value = -5261.0003
$ws.Range("N132").Value = -30333.638

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 3227.5
$ws.Range("I102").Value = 2327.8235
$ws.Range("K102").Value = 2327.8235
$ws.Range("M102").Value = -705.8235

# Row 124
$ws.Range("H124").Value = 41827.145
$ws.Range("J124").Value = 41827.145
$ws.Range("L124").Value = 41827.145
$ws.Range("N124").Value = -51647.145

# Row 126
$ws.Range("H126").Value = 3375.39
$ws.Range("I126").Value = 2849.1758
$ws.Range("J126").Value = 4873.077
$ws.Range("K126").Value = 8547.527399999999
$ws.Range("L126").Value = 14619.231
$ws.Range("M126").Value = -6077.527399999999
$ws.Range("N126").Value = -19559.231

# Row 131
$ws.Range("H131").Value = 31000
$ws.Range("J131").Value = 31000
$ws.Range("L131").Value = 31000
$ws.Range("N131").Value = -41080

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1134.6428
$ws.Range("I16").Value = 1129.6154
$ws.Range("J16").Value = 1200
$ws.Range("K16").Value = 1129.6154
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = -959.6153999999999
$ws.Range("N16").Value = -1540

# Row 26
$ws.Range("H26").Value = 16336
$ws.Range("I26").Value = 5009
$ws.Range("J26").Value = 21999.5
$ws.Range("K26").Value = 5009
$ws.Range("L26").Value = 21999.5
$ws.Range("M26").Value = -4714
$ws.Range("N26").Value = -22589.5

# Row 53
$ws.Range("H53").Value = 10014.667
$ws.Range("I53").Value = 8045
$ws.Range("J53").Value = 10999.5
$ws.Range("K53").Value = 8045
$ws.Range("L53").Value = 10999.5
$ws.Range("M53").Value = -7527
$ws.Range("N53").Value = -12035.5

# Row 104
$ws.Range("H104").Value = 25435
$ws.Range("J104").Value = 25435
$ws.Range("L104").Value = 25435
$ws.Range("N104").Value = -32423

# Row 122
$ws.Range("H122").Value = 5033.2666
$ws.Range("I122").Value = 4250
$ws.Range("J122").Value = 8166.3335
$ws.Range("K122").Value = 12750
$ws.Range("L122").Value = 24499.0005
$ws.Range("M122").Value = -10300
$ws.Range("N122").Value = -29399.0005

# Row 127
$ws.Range("H127").Value = 26456.25
$ws.Range("J127").Value = 26456.25
$ws.Range("L127").Value = 26456.25
$ws.Range("N127").Value = -36376.25

# Row 132
$ws.Range("H132").Value = 7828.8
$ws.Range("I132").Value = 3037.8
$ws.Range("K132").Value = 9113.400000000001
$ws.Range("M132").Value = -6583.400000000001

# Row 136
$ws.Range("H136").Value = 4610.143
$ws.Range("I136").Value = 1334.7778
$ws.Range("K136").Value = 4004.3334
$ws.Range("M136").Value = -1454.3334

$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 333342180
$ws.Range("J5").Value = 13251
$ws.Range("L5").Value = 13251
$ws.Range("N5").Value = -13475

# Row 24
$ws.Range("H24").Value = 9866.666999999999
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 9866.666999999999
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 9866.666999999999
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -10326.667

# Row 42
$ws.Range("H42").Value = 30000
$ws.Range("J42").Value = 30000
$ws.Range("L42").Value = 30000
$ws.Range("N42").Value = -30756

# Row 94
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

# Row 125
$ws.Range("H125").Value = 40482.5
$ws.Range("J125").Value = 40482.5
$ws.Range("L125").Value = 40482.5
$ws.Range("N125").Value = -50322.5

# Row 126
$ws.Range("H126").Value = 2171.1538
$ws.Range("I126").Value = 1120.2273
$ws.Range("J126").Value = 7951.25
$ws.Range("K126").Value = 3360.6819
$ws.Range("L126").Value = 23853.75
$ws.Range("M126").Value = -890.6819
$ws.Range("N126").Value = -28793.75

# Row 136
$ws.Range("H136").Value = 6936.5312
$ws.Range("I136").Value = 7439.1113
$ws.Range("J136").Value = 6290.357
$ws.Range("K136").Value = 22317.3339
$ws.Range("L136").Value = 18871.071
$ws.Range("M136").Value = -19767.3339
$ws.Range("N136").Value = -23971.071
